$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B18: "151 à 174" -> "152 à 175"
$ws.Range("B18").Value = "152 à 175"

# Insert a new blank row before row 26, pushing rows 26-32 down to 27-33
$ws.Rows("26").Insert() | Out-Null

# Move the selection to A19, matching the post-edit cursor position
$ws.Range("A19").Select() | Out-Null
